$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("A2").Value = 111813153
$ws.Range("B2").Value = 56398
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("J2").ClearContents()
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").Value = "färska spår"
$ws.Range("AF2").ClearContents()

# --- Row 3 updates ---
$ws.Range("A3").Value = 111813166
$ws.Range("B3").Value = 90087
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 3298
$ws.Range("F3").Value = "Trådticka"
$ws.Range("G3").Value = "Climacocystis borealis"
$ws.Range("H3").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("J3").Style = "Normal"
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("AF3").Style = "Normal"
